# Updated cryptos list on Wed Jul  3 11:00:20 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# re-sorts the Monero/Aptos (rows 38-39) and Cosmos/FirstDigitalUSD
# (rows 49-50) pairs to their new ranking order.
#
# Price values are stored as free-form text (e.g. "27.96", "3.342.01"), not
# numbers, so for any new price string that Excel would otherwise silently
# re-parse as a numeric literal (dropping trailing zeros / significant
# digits), the cell is switched to the Text number format *before* the
# value is assigned. That keeps the literal characters intact instead of
# Excel coercing them into a float.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.283.23'
$ws.Range('E2').Value = '  -3.86%  '
$ws.Range('D3').Value = '3.311.18'
$ws.Range('E3').Value = '  -4.04%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '560.70'
$ws.Range('E5').Value = '  -3.15%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.43'
$ws.Range('E6').Value = '  -2.90%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.312.51'
$ws.Range('E8').Value = '  -4.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.483'
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('E10').Value = '  -2.63%  '
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('E12').Value = '  -1.20%  '
$ws.Range('D13').Value = '3.876.50'
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '27.29'
$ws.Range('E15').Value = '  -3.76%  '
$ws.Range('D16').Value = '3.315.33'
$ws.Range('E16').Value = '  -3.95%  '
$ws.Range('E17').Value = '  -2.95%  '
$ws.Range('D18').Value = '60.243.84'
$ws.Range('E18').Value = '  -3.97%  '
$ws.Range('E19').Value = '  -3.73%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.37'
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '8.68'
$ws.Range('E21').Value = '  -3.59%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '375.08'
$ws.Range('E22').Value = '  -3.10%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '74.27'
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('E24').Value = '  -2.57%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').Value = '3.476.98'
$ws.Range('E26').Value = '  -2.97%  '
$ws.Range('E27').Value = '  -7.86%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.173'
$ws.Range('E28').Value = '  -5.31%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('E30').Value = '  -5.34%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  -2.79%  '
$ws.Range('E33').Value = '  -4.63%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '22.65'
$ws.Range('E34').Value = '  -2.58%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.29'
$ws.Range('E35').Value = '  -3.42%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.19'
$ws.Range('E36').Value = '  -3.11%  '
$ws.Range('E37').Value = '  -4.98%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.79'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '166.52'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '27.96'
$ws.Range('E40').Value = '  -13.20%  '
$ws.Range('D41').Value = '3.342.01'
$ws.Range('E41').Value = '  -3.99%  '
$ws.Range('E42').Value = '  -4.75%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '41.96'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.754'
$ws.Range('E44').Value = '  -4.05%  '
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('E47').Value = '  -4.92%  '
$ws.Range('D48').Value = '2.393.86'
$ws.Range('E48').Value = '  -6.99%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.61'
$ws.Range('E50').Value = '  -4.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '21.79'
$ws.Range('E51').Value = '  -3.49%  '
